$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: restyle existing row 18 (closing row of the previous block) to the
# bordered "closing" style, matching rows 7/12/14/16 ---
$ws.Range("A14:E14").Copy()
$ws.Range("A18:E18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Step 2: apply formatting templates (copy formats only) for the new rows ---
# Rows 19/22/25/26/27: "header" style - A,B,C,D,E unbordered, like row 17
$ws.Range("A17:E17").Copy()
$ws.Range("A19:E19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A22:E22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A25:E25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A26:E26").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A27:E27").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Rows 20/23: "middle" style - B,C,D,E only, unbordered, like row 11
$ws.Range("B11:E11").Copy()
$ws.Range("B20:E20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B23:E23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Rows 21/24: "closing" style - A,B,C,D,E bordered, like the newly restyled row 18
$ws.Range("A18:E18").Copy()
$ws.Range("A21:E21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A24:E24").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row 27 only uses column A; clear the copied formatting from B:E so those cells
# do not exist at all
$ws.Range("B27:E27").Clear()

$excel.CutCopyMode = $false

# --- Step 3: set the cell values/text for the new rows ---

# Row 19
$ws.Range("A19").Value = 'SCRIPT/P01P04A/um1201.ssb'
$ws.Range("B19").Value = 301
$ws.Range("C19").Value = ' I\''m really enjoying the \"Big\nTreasure\" [CS:I]Prize Tickets[CR] at the Recycle Shop.\nI can\''t get enough!'
$ws.Range("D19").Value = ' Мне нравится идея [CS:I]Призовых\nБилетов[CR] Центра Переработки. Я только и\nделаю, что собираю их!'
$ws.Range("E19").Value = ' Íîå îñàâéóòÿ éäåÿ [CS:I]Ðñéèïâúö\nÁéìåóïâ[CR] Øåîóñà Ðåñåñàáïóëé. Ÿ óïìûëï é\näåìàý, œóï òïáéñàý éö!'

# Row 20
$ws.Range("B20").Value = 308
$ws.Range("C20").Value = ' I wonder if I\''ll win on my very\nfirst [CS:I]Prize Ticket[CR]. Ho-ho-ho!'
$ws.Range("D20").Value = ' Интересно, что я получу за свой\nпервый [CS:I]Призовой Билет[CR]. Хо-хо-хо!'
$ws.Range("E20").Value = ' Éîóåñåòîï, œóï ÿ ðïìôœô èà òâïê\nðåñâúê [CS:I]Ðñéèïâïê Áéìåó[CR]. Öï-öï-öï!'

# Row 21
$ws.Range("B21").Value = 315
$ws.Range("C21").Value = ' One day, I\''ll win big!'
$ws.Range("D21").Value = ' Однажды я выиграю по-крупному!'
$ws.Range("E21").Value = ' Ïäîàçäú ÿ âúéãñàý ðï-ëñôðîïíô!'

# Row 22
$ws.Range("A22").Value = 'SCRIPT/P01P04A/um1301.ssb'
$ws.Range("B22").Value = 301
$ws.Range("C22").Value = ' I\''m really enjoying the \"Big\nTreasure\" [CS:I]Prize Tickets[CR] at the Recycle Shop.\nI can\''t get enough!'
$ws.Range("D22").Value = ' Мне нравится идея [CS:I]Призовых\nБилетов[CR] Центра Переработки. Я только и\nделаю, что собираю их!'
$ws.Range("E22").Value = ' Íîå îñàâéóòÿ éäåÿ [CS:I]Ðñéèïâúö\nÁéìåóïâ[CR] Øåîóñà Ðåñåñàáïóëé. Ÿ óïìûëï é\näåìàý, œóï òïáéñàý éö!'

# Row 23
$ws.Range("B23").Value = 308
$ws.Range("C23").Value = ' I wonder if I\''ll win on my very\nfirst [CS:I]Prize Ticket[CR]. Ho-ho-ho!'
$ws.Range("D23").Value = ' Интересно, что я получу за свой\nпервый [CS:I]Призовой Билет[CR]. Хо-хо-хо!'
$ws.Range("E23").Value = ' Éîóåñåòîï, œóï ÿ ðïìôœô èà òâïê\nðåñâúê [CS:I]Ðñéèïâïê Áéìåó[CR]. Öï-öï-öï!'

# Row 24
$ws.Range("B24").Value = 315
$ws.Range("C24").Value = ' One day, I\''ll win big!'
$ws.Range("D24").Value = ' Однажды я выиграю по-крупному!'
$ws.Range("E24").Value = ' Ïäîàçäú ÿ âúéãñàý ðï-ëñôðîïíô!'

# Row 25
$ws.Range("A25").Value = 'SCRIPT/P01P04A/um1402.ssb'
$ws.Range("B25").Value = 279
$ws.Range("C25").Value = ' I love to hear [CS:N]Wobbuffet[CR]\nblurt out, \"That\''s right!\" when someone wins\non a [CS:I]Prize Ticket[CR].'
$ws.Range("D25").Value = ' Я обожаю, когда [CS:N]Воббаффет[CR]\nвыкрикивает: \"Всё так!\", когда кто-нибудь\nполучает выигрыш за [CS:I]Призовой Билет[CR].'
$ws.Range("E25").Value = ' Ÿ ïáïçàý, ëïãäà [CS:N]Âïááàõõåó[CR]\nâúëñéëéâàåó: \"Âòæ óàë!\", ëïãäà ëóï-îéáôäû\nðïìôœàåó âúéãñúš èà [CS:I]Ðñéèïâïê Áéìåó[CR].'

# Row 26
$ws.Range("A26").Value = 'SCRIPT/P01P04A/um1502.ssb'
$ws.Range("B26").Value = 282
$ws.Range("C26").Value = ' Oh... I\''d love to hear her yell\nthat again…'
$ws.Range("D26").Value = ' О... Я так хочу снова услышать\nеё крик...'
$ws.Range("E26").Value = ' Ï... Ÿ óàë öïœô òîïâà ôòìúšàóû\nåæ ëñéë...'

# Row 27
$ws.Range("A27").Value = 'SCRIPT/P01P04A/um1602.ssb'


# --- Step 4: set explicit row heights to match the target layout ---
$ws.Rows.Item(19).RowHeight = 43.2
$ws.Rows.Item(20).RowHeight = 31.8
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 31.8
$ws.Rows.Item(25).RowHeight = 52.2
$ws.Rows.Item(26).RowHeight = 43.2
$ws.Rows.Item(27).RowHeight = 43.2

# --- Step 5: update sheet view (scroll position + selected cell) ---
$ws.Range("E26").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
